$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "27.464.44"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.566.41"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D5").Value = "'208.31"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'0.0591"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.789.53"
$ws.Range("D13").Value = "1.558.42"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "'63.59"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "27.458.77"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "'213.28"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "'9.55"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'153.15"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'0.0469"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "1.375.46"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").Value = "'0.954"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D42").Value = "'0.977"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").Value = "'64.13"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "'2.17"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "1.702.00"
$ws.Range("D48").Value = "'85.50"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "'0.0958"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "'0.0495"
$ws.Range("E51").Value = "  -0.68%  "
